# Import Aspirants functionality commit & Removed import from Employees section
#
# Rebuilds the user-import template: replaces the 5-column "Employee" style
# header/sample row with the 18-column "Aspirant" layout (FirstName..Team,
# including the duplicated CurrentAddress* block, Note/Interest columns),
# a sample data row, cell alignment, a mailto hyperlink on the sample email,
# and dropdown (list) data validations for State and Team.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1) — written in this specific order so the workbook's
#    shared-string table is built up the same way it was originally authored
#    (base fields first, then the "Current Address" duplicate block + notes).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "MiddleName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Address 1"
$ws.Range("F1").Value = "Address 2"
$ws.Range("G1").Value = "City"
$ws.Range("H1").Value = "State"
$ws.Range("I1").Value = "Pincode"
$ws.Range("O1").Value = "Mobile"
$ws.Range("R1").Value = "Team"

# ---------------------------------------------------------------------------
# 2) Sample data row (row 2) — core identity + first address block.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "Nagpur"
$ws.Range("F2").Value = "MUMBAI"
$ws.Range("M2").Value = "Maharashtra"
$ws.Range("A2").Value = "Ramakant"
$ws.Range("B2").Value = "Shyam"
$ws.Range("C2").Value = "Chandel"
$ws.Range("D2").Value = "ramakant.chandel@perficient.com"

# ---------------------------------------------------------------------------
# 3) Remaining header cells — the "Current Address" duplicate block + notes.
# ---------------------------------------------------------------------------
$ws.Range("J1").Value = "CurrentAddress1"
$ws.Range("K1").Value = "CurrentAddress2"
$ws.Range("L1").Value = "CurrentCity"
$ws.Range("M1").Value = "CurrentStateId"
$ws.Range("N1").Value = "CurrentPincode"
$ws.Range("P1").Value = "Note"
$ws.Range("Q1").Value = "Interest"

# ---------------------------------------------------------------------------
# 4) Remaining sample data cells.
# ---------------------------------------------------------------------------
$ws.Range("P2").Value = "This is a note"
$ws.Range("Q2").Value = "HTML/CSS, Design"
$ws.Range("H2").Value = "Madhya Pradesh"
$ws.Range("R2").Value = "Magento"

$ws.Range("G2").Value = "Nagpur"
$ws.Range("J2").Value = "Nagpur"
$ws.Range("L2").Value = "Nagpur"
$ws.Range("K2").Value = "MUMBAI"
$ws.Range("I2").Value = 440024
$ws.Range("N2").Value = 440024
$ws.Range("O2").Value = 9960160804

# ---------------------------------------------------------------------------
# 5) Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.59
$ws.Columns.Item(2).ColumnWidth = 12.25
$ws.Columns.Item(3).ColumnWidth = 12.09
$ws.Columns.Item(4).ColumnWidth = 32.76
$ws.Columns.Item(5).ColumnWidth = 12.42
$ws.Columns.Item(6).ColumnWidth = 11.42
$ws.Columns.Item(7).ColumnWidth = 9.92
$ws.Columns.Item(8).ColumnWidth = 18.42
$ws.Columns.Item(10).ColumnWidth = 15.76
$ws.Columns.Item(11).ColumnWidth = 17.59
$ws.Columns.Item(12).ColumnWidth = 17.42
$ws.Columns.Item(13).ColumnWidth = 17.59
$ws.Columns.Item(14).ColumnWidth = 18.59
$ws.Columns.Item(15).ColumnWidth = 12.92
$ws.Columns.Item(16).ColumnWidth = 12.92
$ws.Columns.Item(17).ColumnWidth = 16.09

# ---------------------------------------------------------------------------
# 6) Alignment — every used cell (A1:R2) is left/top aligned. Build the combo
#    via a throwaway named style so only ONE new cellXf is created instead of
#    one per property assignment.
# ---------------------------------------------------------------------------
$alignStyle = $wb.Styles.Add("LeftTopAlign")
$alignStyle.HorizontalAlignment = -4131
$alignStyle.VerticalAlignment = -4160
$ws.Range("A1:R2").Style = "LeftTopAlign"

# ---------------------------------------------------------------------------
# 7) Hyperlink on the sample email cell, then restore its left/top alignment
#    (adding the hyperlink assigns the built-in Hyperlink font/style, which
#    would otherwise clobber the alignment applied above).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:ramakant.chandel@perficient.com")
$ws.Range("D2").HorizontalAlignment = -4131
$ws.Range("D2").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 8) Data validation dropdowns.
# ---------------------------------------------------------------------------
$ws.Range("H2").Validation.Add(3, 1, 1, """Maharashta, Delhi, Madhya Pradesh, Tamil Nadu""")
$ws.Range("M2").Validation.Add(3, 1, 1, """Maharashta, Delhi, Madhya Pradesh, Tamil Nadu""")
$ws.Range("R2").Validation.Add(3, 1, 1, """AEM, .Net Technology, Magento, UI team, ALL, PHP Team""")

# ---------------------------------------------------------------------------
# 9) Page setup + selection.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("S16").Select()

Write-Host "Import Aspirants template applied"
